$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "mapsto"

# Update the "mapsto" labels in column B with clearer/renamed stage names
$ws.Range("B2").Value = "wake"
$ws.Range("B3").Value = "rem"
$ws.Range("B4").Value = "stage1"
$ws.Range("B5").Value = "stage2"
$ws.Range("B6").Value = "sws"
$ws.Range("B7").Value = "movement"
$ws.Range("B8").Value = "unknown"
$ws.Range("B9").Value = "artifact"

# Remove the old numeric "mapsto" codes in column C entirely
$ws.Range("C1:C9").Clear()

# Restore the selection Excel had after the edit
$ws.Range("B7").Select()
